$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.122.51"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "3.139.32"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.77"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.78"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D8").Value = "3.139.26"
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.468"
$ws.Range("E9").Value = "  +4.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.28"
$ws.Range("E10").Value = "  +1.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.107"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  +4.36%  "
$ws.Range("D13").Value = "3.679.21"
$ws.Range("E14").Value = "  +1.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.51"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000164"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "58.188.52"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").Value = "3.143.17"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.04"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.68"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.15"
$ws.Range("E21").Value = "  +2.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "359.71"
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.98"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.506"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  -1.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "0.0₃0876"
$ws.Range("E28").Value = "  -4.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.31"
$ws.Range("E29").Value = "  -2.47%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.88"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.09"
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.58"
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.03"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("E34").Value = "  -2.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.39"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.08"
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.87"
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.27"
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.69"
$ws.Range("E39").Value = "  +4.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0671"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "2.503.56"
$ws.Range("E41").Value = "  +8.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.704"
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.00"
$ws.Range("E43").Value = "  -4.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "37.41"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "3.182.98"
$ws.Range("E45").Value = "  +0.28%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0269"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.988"
$ws.Range("E48").Value = "  +2.53%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.78"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.737"
$ws.Range("E51").Value = "  -4.49%  "
